# "Lagt in tider jag kan" - fill in availability (":D" = available, ":(" = not
# available) for the "Mµ" row of both weeks on the schema sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 1 - row 12 ("Mµ")
$ws.Range("B12").Value = ":D"
$ws.Range("C12").Value = ":D"
$ws.Range("D12").Value = ":D"
$ws.Range("E12").Value = ":("
$ws.Range("F12").Value = ":("
$ws.Range("G12").Value = ":("
$ws.Range("H12").Value = ":D"
$ws.Range("I12").Value = ":("
$ws.Range("J12").Value = ":D"
$ws.Range("K12").Value = ":D"

# Week 2 - row 20 ("Mµ")
$ws.Range("B20").Value = ":D"
$ws.Range("C20").Value = ":D"
$ws.Range("D20").Value = ":("
$ws.Range("E20").Value = ":("
$ws.Range("F20").Value = ":D"
$ws.Range("G20").Value = ":D"
$ws.Range("H20").Value = ":("
$ws.Range("I20").Value = ":("
$ws.Range("J20").Value = ":D"
$ws.Range("K20").Value = ":D"

# Leave the active selection on the last cell that was filled in.
[void]$ws.Range("B20").Select()
